$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Apply the formatting once (bold font, thin box border, center/top alignment)
# to B1, then copy/paste the resulting format onto A2 so both cells end up
# sharing a single style entry instead of each accreting its own xf.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Borders.LineStyle = 1
$b1.Borders.Weight = 2
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160

$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
